$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.43616840798465
$ws.Range("D2").Value = 3.431135433691324
$ws.Range("E2").Value = 19.75547178734143
$ws.Range("F2").Value = 25.97006422263257
$ws.Range("G2").Value = 34.5223991359973
$ws.Range("H2").Value = 13.85183328071549
$ws.Range("I2").Value = 22.55805929405394
$ws.Range("L2").Value = 10.1090585160351
$ws.Range("N2").Value = 17.68618808484334
$ws.Range("B3").Value = 15.09658988651523
$ws.Range("D3").Value = 3.450456682539366
$ws.Range("E3").Value = 19.34131854567429
$ws.Range("F3").Value = 25.33872057111664
$ws.Range("G3").Value = 33.1945802744435
$ws.Range("H3").Value = 13.74167325498498
$ws.Range("I3").Value = 22.72911852487826
$ws.Range("L3").Value = 9.91086866741111
$ws.Range("N3").Value = 17.70715894152752
$ws.Range("B4").Value = 14.88655162169721
$ws.Range("D4").Value = 3.463203458708285
$ws.Range("E4").Value = 19.08070159868165
$ws.Range("F4").Value = 24.95362956415997
$ws.Range("G4").Value = 32.3664761644595
$ws.Range("H4").Value = 13.67826265619505
$ws.Range("I4").Value = 22.83931618153206
$ws.Range("L4").Value = 9.789233909878485
$ws.Range("N4").Value = 17.72220381421018
$ws.Range("B5").Value = 14.80069245059844
$ws.Range("D5").Value = 3.468618912312218
$ws.Range("E5").Value = 18.97299722032744
$ws.Range("F5").Value = 24.79761184677212
$ws.Range("G5").Value = 32.02642329200862
$ws.Range("H5").Value = 13.65350835224798
$ws.Range("I5").Value = 22.88552517306984
$ws.Range("L5").Value = 9.739748252531157
$ws.Range("N5").Value = 17.728881419052
$ws.Range("B6").Value = 14.78642304597078
$ws.Range("D6").Value = 3.469531462552303
$ws.Range("E6").Value = 18.95502504998639
$ws.Range("F6").Value = 24.77176808777366
$ws.Range("G6").Value = 31.96982062141986
$ws.Range("H6").Value = 13.64946413351279
$ws.Range("I6").Value = 22.89327691127783
$ws.Range("L6").Value = 9.731538121196715
$ws.Range("N6").Value = 17.73002328848119
$ws.Range("B7").Value = 14.88539461496209
$ws.Range("D7").Value = 3.46327559990772
$ws.Range("E7").Value = 19.07925501474057
$ws.Range("F7").Value = 24.95152141360915
$ws.Range("G7").Value = 32.36189972093545
$ws.Range("H7").Value = 13.67792438541635
$ws.Range("I7").Value = 22.83993409375164
$ws.Range("L7").Value = 9.788566108701307
$ws.Range("N7").Value = 17.72229165541897
$ws.Range("B8").Value = 15.31947073181243
$ws.Range("D8").Value = 3.437613362199095
$ws.Range("E8").Value = 19.6140369165083
$ws.Range("F8").Value = 25.75200480352845
$ws.Range("G8").Value = 34.06758754447554
$ws.Range("H8").Value = 13.81298697155228
$ws.Range("I8").Value = 22.61597022135967
$ws.Range("L8").Value = 10.0407515792684
$ws.Range("N8").Value = 17.69296955543154
$ws.Range("B9").Value = 16.15349406636456
$ws.Range("D9").Value = 3.394350692933696
$ws.Range("E9").Value = 20.60885693586374
$ws.Range("F9").Value = 27.33064474693009
$ws.Range("G9").Value = 37.28599433533717
$ws.Range("H9").Value = 14.11027946486744
$ws.Range("I9").Value = 22.2176169844632
$ws.Range("L9").Value = 10.53281766701134
$ws.Range("N9").Value = 17.65261608952065
$ws.Range("B10").Value = 16.74942247961569
$ws.Range("D10").Value = 3.366942293377924
$ws.Range("E10").Value = 21.30230063082425
$ws.Range("F10").Value = 28.4813291225588
$ws.Range("G10").Value = 39.54439298214047
$ws.Range("H10").Value = 14.34684265973221
$ws.Range("I10").Value = 21.9496238512246
$ws.Range("L10").Value = 10.88909016148714
$ws.Range("N10").Value = 17.63333970058979
$ws.Range("B11").Value = 17.01564733780085
$ws.Range("D11").Value = 3.355440183023747
$ws.Range("E11").Value = 21.60876130775726
$ws.Range("F11").Value = 28.99998478710416
$ws.Range("G11").Value = 40.5437022988596
$ws.Range("H11").Value = 14.45801835678363
$ws.Range("I11").Value = 21.83301888191634
$ws.Range("L11").Value = 11.04928661365063
$ws.Range("N11").Value = 17.62680406208001
$ws.Range("B12").Value = 17.11566666640588
$ws.Range("D12").Value = 3.351224987900948
$ws.Range("E12").Value = 21.72345089604973
$ws.Range("F12").Value = 29.19547998339627
$ws.Range("G12").Value = 40.91772985704215
$ws.Range("H12").Value = 14.50059585482345
$ws.Range("I12").Value = 21.78962332244282
$ws.Range("L12").Value = 11.1096216902291
$ws.Range("N12").Value = 17.62464868661502
$ws.Range("B13").Value = 17.09416247279743
$ws.Range("D13").Value = 3.352126537028883
$ws.Range("E13").Value = 21.69881196878719
$ws.Range("F13").Value = 29.15342031488887
$ws.Range("G13").Value = 40.83737671868437
$ws.Range("H13").Value = 14.49140534371157
$ws.Range("I13").Value = 21.79893556671168
$ws.Range("L13").Value = 11.09664292703962
$ws.Range("N13").Value = 17.62509870045134
$ws.Range("B14").Value = 17.0238924237501
$ws.Range("D14").Value = 3.355090573840706
$ws.Range("E14").Value = 21.61822452324602
$ws.Range("F14").Value = 29.01608769853098
$ws.Range("G14").Value = 40.57456361178767
$ws.Range("H14").Value = 14.46151186097293
$ws.Range("I14").Value = 21.82943348345705
$ws.Range("L14").Value = 11.05425730977767
$ws.Range("N14").Value = 17.62662034365638
$ws.Range("B15").Value = 16.98074386383981
$ws.Range("D15").Value = 3.356924461432359
$ws.Range("E15").Value = 21.56868329933647
$ws.Range("F15").Value = 28.93184297367119
$ws.Range("G15").Value = 40.41300161969521
$ws.Range("H15").Value = 14.44326240600667
$ws.Range("I15").Value = 21.84821325653193
$ws.Range("L15").Value = 11.02825052432028
$ws.Range("N15").Value = 17.62759395727847
$ws.Range("B16").Value = 16.73191791056711
$ws.Range("D16").Value = 3.367713549342901
$ws.Range("E16").Value = 21.28208589953143
$ws.Range("F16").Value = 28.44731818670007
$ws.Range("G16").Value = 39.47848987810815
$ws.Range("H16").Value = 14.33964581269072
$ws.Range("I16").Value = 21.95735079460558
$ws.Range("L16").Value = 10.87857817402322
$ws.Range("N16").Value = 17.63381162260136
$ws.Range("B17").Value = 16.5779555858718
$ws.Range("D17").Value = 3.374580787222003
$ws.Range("E17").Value = 21.10391530394166
$ws.Range("F17").Value = 28.14868890819871
$ws.Range("G17").Value = 38.89774523709823
$ws.Range("H17").Value = 14.27696844143351
$ws.Range("I17").Value = 22.02566020211871
$ws.Range("L17").Value = 10.78623652684064
$ws.Range("N17").Value = 17.63819682584853
$ws.Range("B18").Value = 16.48894775299913
$ws.Range("D18").Value = 3.378621484291104
$ws.Range("E18").Value = 21.00059379324559
$ws.Range("F18").Value = 27.97648589824849
$ws.Range("G18").Value = 38.5610944709475
$ws.Range("H18").Value = 14.24125547840439
$ws.Range("I18").Value = 22.06544960170908
$ws.Range("L18").Value = 10.73295098045334
$ws.Range("N18").Value = 17.64092938649819
$ws.Range("B19").Value = 16.45873645492505
$ws.Range("D19").Value = 3.380005152977028
$ws.Range("E19").Value = 20.96546834953803
$ws.Range("F19").Value = 27.91811237815082
$ws.Range("G19").Value = 38.44667162654583
$ws.Range("H19").Value = 14.22922266739476
$ws.Range("I19").Value = 22.07900750927088
$ws.Range("L19").Value = 10.71488151870298
$ws.Range("N19").Value = 17.64189075341006
$ws.Range("B20").Value = 16.59439270469071
$ws.Range("D20").Value = 3.373840345500845
$ws.Range("E20").Value = 21.12296956917821
$ws.Range("F20").Value = 28.18052560935083
$ws.Range("G20").Value = 38.95984056982266
$ws.Range("H20").Value = 14.2836058816399
$ws.Range("I20").Value = 22.01833685958912
$ws.Range("L20").Value = 10.79608479897159
$ws.Range("N20").Value = 17.63770825895406
$ws.Range("B21").Value = 17.0445547163681
$ws.Range("D21").Value = 3.354216141831964
$ws.Range("E21").Value = 21.64193244451963
$ws.Range("F21").Value = 29.05645193765803
$ws.Range("G21").Value = 40.6518799211464
$ws.Range("H21").Value = 14.47027961734077
$ws.Range("I21").Value = 21.82045489650561
$ws.Range("L21").Value = 11.06671633938293
$ws.Range("N21").Value = 17.62616474138252
$ws.Range("B22").Value = 17.33409300415317
$ws.Range("D22").Value = 3.342209692790002
$ws.Range("E22").Value = 21.97314732582679
$ws.Range("F22").Value = 29.62353493926058
$ws.Range("G22").Value = 41.73200856762542
$ws.Range("H22").Value = 14.59504932576633
$ws.Range("I22").Value = 21.69555741122236
$ws.Range("L22").Value = 11.24165408274882
$ws.Range("N22").Value = 17.62048205607414
$ws.Range("B23").Value = 17.18001765032039
$ws.Range("D23").Value = 3.348542312363586
$ws.Range("E23").Value = 21.79712083774185
$ws.Range("F23").Value = 29.32143305048174
$ws.Range("G23").Value = 41.1579825131189
$ws.Range("H23").Value = 14.52821578181418
$ws.Range("I23").Value = 21.76181314675111
$ws.Range("L23").Value = 11.14848204660362
$ws.Range("N23").Value = 17.62334521712881
$ws.Range("B24").Value = 16.58696301112007
$ws.Range("D24").Value = 3.37417481079071
$ws.Range("E24").Value = 21.11435789928612
$ws.Range("F24").Value = 28.16613383051856
$ws.Range("G24").Value = 38.93177587947242
$ws.Range("H24").Value = 14.28060409218398
$ws.Range("I24").Value = 22.02164613091004
$ws.Range("L24").Value = 10.79163300460626
$ws.Range("N24").Value = 17.63792848130918
$ws.Range("B25").Value = 15.93038730209424
$ws.Range("D25").Value = 3.405291174817218
$ws.Range("E25").Value = 20.34599500619175
$ws.Range("F25").Value = 26.90417764918251
$ws.Range("G25").Value = 36.43208710321555
$ws.Range("H25").Value = 14.02654278326983
$ws.Range("I25").Value = 22.32103287790751
$ws.Range("L25").Value = 10.40036099726935
$ws.Range("N25").Value = 17.66170532845616
